$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-09-28 Sunday" "2025-09-29 Monday"

Replace-Text "71×29=2059" "77×53=4081"
Replace-Text "59×24=1416" "41×27=1107"
Replace-Text "33×92=3036" "57×20=1140"
Replace-Text "41×21=861" "74×52=3848"
Replace-Text "52×59=3068" "93×17=1581"

Replace-Text "41×16=656" "63×73=4599"
Replace-Text "91×18=1638" "52×51=2652"
Replace-Text "61×72=4392" "68×71=4828"
Replace-Text "50×65=3250" "18×85=1530"
Replace-Text "96×33=3168" "45×85=3825"

Replace-Text "55×25=1375" "23×69=1587"
Replace-Text "93×36=3348" "78×72=5616"
Replace-Text "54×66=3564" "30×22=660"
Replace-Text "98×11=1078" "55×29=1595"
Replace-Text "18×52=936" "22×70=1540"

Replace-Text "77×96=7392" "60×80=4800"
Replace-Text "36×63=2268" "20×99=1980"
Replace-Text "26×74=1924" "36×72=2592"
Replace-Text "17×93=1581" "86×50=4300"
Replace-Text "19×72=1368" "90×83=7470"

Replace-Text "30×24=720" "69×89=6141"
Replace-Text "97×83=8051" "37×72=2664"
Replace-Text "40×15=600" "71×80=5680"
Replace-Text "41×22=902" "94×78=7332"
Replace-Text "33×68=2244" "16×18=288"
